$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-04 18:02:22"
$ws.Range("G2").Value = "125 cm"
$ws.Range("E3").Value = "2026-02-04 18:02:24"
$ws.Range("O3").Value = "-4.8 °C"
$ws.Range("E4").Value = "2026-02-04 18:02:27"
$ws.Range("H4").Value = "'82%"
$ws.Range("J4").Value = "993.3 hPa"
$ws.Range("O4").Value = "5.8 °C"
$ws.Range("E5").Value = "2026-02-04 18:02:30"
$ws.Range("H5").Value = "'68%"
$ws.Range("J5").Value = "992.7 hPa"
$ws.Range("E6").Value = "2026-02-04 18:02:33"
$ws.Range("H6").Value = "'66%"
$ws.Range("O6").Value = "11.0 °C"
$ws.Range("E7").Value = "2026-02-04 18:03:35"
$ws.Range("H7").Value = "'83%"
$ws.Range("J7").Value = "993.7 hPa"
$ws.Range("E8").Value = "2026-02-04 18:02:38"
$ws.Range("E9").Value = "2026-02-04 18:02:41"
$ws.Range("E10").Value = "2026-02-04 18:02:44"
$ws.Range("O10").Value = "8.3 °C"
$ws.Range("E11").Value = "2026-02-04 18:02:46"
$ws.Range("J11").Value = "995.2 hPa"
$ws.Range("O11").Value = "0.6 °C"
$ws.Range("E12").Value = "2026-02-04 18:02:49"
$ws.Range("H12").Value = "'80%"
$ws.Range("O12").Value = "8.6 °C"
$ws.Range("E13").Value = "2026-02-04 18:02:52"
$ws.Range("O13").Value = "7.3 °C"
$ws.Range("E14").Value = "2026-02-04 18:02:54"
$ws.Range("E15").Value = "2026-02-04 18:02:57"
$ws.Range("H15").Value = "'82%"
$ws.Range("J15").Value = "993.4 hPa"
$ws.Range("O15").Value = "6.3 °C"
$ws.Range("E16").Value = "2026-02-04 18:03:00"
$ws.Range("E17").Value = "2026-02-04 18:03:03"
$ws.Range("E18").Value = "2026-02-04 18:03:06"
$ws.Range("E19").Value = "2026-02-04 18:03:08"
$ws.Range("O19").Value = "7.0 °C"
$ws.Range("E20").Value = "2026-02-04 18:03:11"
$ws.Range("O20").Value = "-4.6 °C"
$ws.Range("E21").Value = "2026-02-04 18:03:14"
$ws.Range("O21").Value = "6.1 °C"
$ws.Range("E22").Value = "2026-02-04 18:03:17"
$ws.Range("O22").Value = "8.4 °C"
$ws.Range("E23").Value = "2026-02-04 18:03:19"
$ws.Range("E24").Value = "2026-02-04 18:03:22"
$ws.Range("K24").Value = "9.4 MJ/m2"
$ws.Range("E25").Value = "2026-02-04 18:03:25"
$ws.Range("E26").Value = "2026-02-04 18:03:27"
$ws.Range("H26").Value = "'68%"
$ws.Range("O26").Value = "-2.0 °C"
$ws.Range("E27").Value = "2026-02-04 18:03:30"
$ws.Range("J27").Value = "993.0 hPa"
$ws.Range("O27").Value = "10.8 °C"
$ws.Range("E28").Value = "2026-02-04 18:03:33"
$ws.Range("H28").Value = "'84%"
$ws.Range("J28").Value = "994.8 hPa"
$ws.Range("O28").Value = "2.4 °C"
$ws.Range("E29").Value = "2026-02-04 18:03:36"
$ws.Range("H29").Value = "'76%"
$ws.Range("O29").Value = "7.2 °C"
$ws.Range("E30").Value = "2026-02-04 18:03:38"
$ws.Range("H30").Value = "'76%"
$ws.Range("E31").Value = "2026-02-04 18:03:41"
$ws.Range("J31").Value = "994.7 hPa"
$ws.Range("O31").Value = "4.2 °C"
$ws.Range("E32").Value = "2026-02-04 18:03:44"
$ws.Range("O32").Value = "10.4 °C"
$ws.Range("E33").Value = "2026-02-04 18:03:47"
$ws.Range("E34").Value = "2026-02-04 18:03:49"
$ws.Range("H34").Value = "'88%"
$ws.Range("O34").Value = "3.4 °C"
$ws.Range("E35").Value = "2026-02-04 18:03:52"
$ws.Range("E36").Value = "2026-02-04 18:03:55"
$ws.Range("H36").Value = "'85%"
$ws.Range("O36").Value = "7.1 °C"
